$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "Daily Temperatures" solution-ideas note (row 3, col E) ---
$ws.Range("E3").Value = "- Brute force`n- Two loops, one for day, one for check temperature`n--> (TLE)`n- Monotonic stack`n- Iterate over the input array`n- keep the index and value of value which small than its previous value`n- while the value bigger than stack.peek(), pop the stack and assign index to answer array.`n--> O(2*N) An easier way to think about this is that in the worst case, every element will be pushed and popped once. `n- use Array`n- "

# --- Remove the "Solution component" column (old column F) entirely ---
$ws.Columns.Item(6).Delete()

# --- Resize remaining columns E (Solution ideas) and F (Note, was G) ---
$ws.Columns.Item(5).ColumnWidth = 93.16666666666667
$ws.Columns.Item(6).ColumnWidth = 16.666666666666668

# --- Adjust row heights now that content/columns changed ---
$ws.Rows.Item(2).RowHeight = 150
$ws.Rows.Item(3).RowHeight = 262.5

# --- Change fill color used by the "Easy" conditional formatting rule ---
$rng = $ws.Range("A1:F1048576")
$fc = $rng.FormatConditions
$ruleEasy = $fc.Item(3)
$ruleEasy.Interior.Color = 5296274

# --- Reset view: scroll back to top-left, keep frozen header pane, select E3 ---
$ws.Range("A1").Select()
$ws.Range("E3").Select()

Write-Host "Edit applied"
